$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.980768084526062
$ws.Range("B1").Value = 3.836013317108154
$ws.Range("C1").Value = 2.767821073532104
$ws.Range("D1").Value = 1.766888737678528
$ws.Range("E1").Value = 1.431453824043274
